$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("network")
$ws.Range("A1").Value = "cols protein1/ rows protein2"

$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "optimization_parameters"

$ws2.Range("A1").Value = "optimization_parameter"
$ws2.Range("B1").Value = "value"
$ws2.Range("A2").Value = "species"
$ws2.Range("B2").Value = "Saccharomyces cerevisiae"
$ws2.Range("A3").Value = "taxon_id"
$ws2.Range("B3").Value = 559292
$ws2.Range("A4").Value = "workbookType"
$ws2.Range("B4").Value = "protein-protein-physical-interaction"
